$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New characters to add: Rey, Admiral Gial Ackbar, Lando Calrissian
# Columns: A=Name, B=Species, C=Gender/Droid Type, D=Birth Year, E=Homeworld,
#          F=First Screen Appearance, G=Wiki Link, H=Image Link

$ws.Range("D16").Value = "15 ABY"
$ws.Range("E16").Value = "Jakku"
$ws.Range("G16").Value = "https://starwars.fandom.com/wiki/Rey_Skywalker"
$ws.Range("A16").Value = "Rey "

$ws.Range("B17").Value = "Mon Calamari"
$ws.Range("E17").Value = "Dac"
$ws.Range("G17").Value = "https://starwars.fandom.com/wiki/Gial_Ackbar/Legends"

$ws.Range("A18").Value = "Lando Calrissian"
$ws.Range("D18").Value = "31 BBY"
$ws.Range("E18").Value = "Socorro"

$ws.Range("F16").Value = "Episode VII: The Force Awakens (2015)"

$ws.Range("G18").Value = "https://starwars.fandom.com/wiki/Lando_Calrissian/Legends"

$ws.Range("H18").Value = "https://cdn.vox-cdn.com/thumbor/Iygx7I0plaTWBYzEuzso_zYEYWg=/1400x1400/filters:format(jpeg)/cdn.vox-cdn.com/uploads/chorus_asset/file/6332771/Lando-Calrissian_a679fe1e.0.jpeg"

$ws.Range("A17").Value = "Admiral Gial Ackbar"

$ws.Range("H16").Value = "https://pyxis.nymag.com/v1/imgs/854/921/9242868f6d9e25cfa955b32c070e3cbb7c-13-the-last-jedi-rey-2.rsquare.w700.jpg"

$ws.Range("H17").Value = "https://pyxis.nymag.com/v1/imgs/5fd/b4c/08e5eaa592aaf4c45ca001b680bb827ae7-13-ackbar.rsquare.w700.jpg"

# Remaining cells (all reuse existing shared strings, order doesn't matter)
$ws.Range("B16").Value = "Human"
$ws.Range("C16").Value = "Female"
$ws.Range("C17").Value = "Male"
$ws.Range("D17").Value = "41 BBY"
$ws.Range("F17").Value = "Episode VI: Return of the Jedi (1983)"
$ws.Range("B18").Value = "Human"
$ws.Range("C18").Value = "Male"
$ws.Range("F18").Value = "Episode V: The Empire Strikes Back (1980)"

# Hyperlinks - added after the values, G column first then H column
$ws.Hyperlinks.Add($ws.Range("G16"), "https://starwars.fandom.com/wiki/Rey_Skywalker") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G18"), "https://starwars.fandom.com/wiki/Lando_Calrissian/Legends") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H16"), "https://pyxis.nymag.com/v1/imgs/854/921/9242868f6d9e25cfa955b32c070e3cbb7c-13-the-last-jedi-rey-2.rsquare.w700.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H17"), "https://pyxis.nymag.com/v1/imgs/5fd/b4c/08e5eaa592aaf4c45ca001b680bb827ae7-13-ackbar.rsquare.w700.jpg") | Out-Null

# Apply the same cell formatting (centered) used by the rest of the table
$ws.Range("A16:H18").HorizontalAlignment = -4108
$ws.Range("A16:H18").VerticalAlignment = -4108

$ws.Range("H22").Select() | Out-Null
